$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("suborder")

$ws.Range("A2").Value = "alb"
$ws.Range("B2").Value = "blanco"
$ws.Range("C2").Value = "presencia de horizonte álbico"

$ws.Range("A3").Value = "antr"
$ws.Range("B3").Value = "antrópico"
$ws.Range("C3").Value = "modificado por el hombre"

$ws.Range("A4").Value = "acu"
$ws.Range("B4").Value = "agua"
$ws.Range("C4").Value = "acracterísticas asociadas con la humedad"

$ws.Range("A5").Value = "ar"
$ws.Range("B5").Value = "arar"
$ws.Range("C5").Value = "horizonte mezclado"

$ws.Range("A6").Value = "arg"
$ws.Range("B6").Value = "arcilla"
$ws.Range("C6").Value = "presencia de horizonte argílico"

$ws.Range("A7").Value = "calc"
$ws.Range("B7").Value = "cal"
$ws.Range("C7").Value = "presencia de horizonte cálcico"

$ws.Range("A8").Value = "camb"
$ws.Range("B8").Value = "cambiar"
$ws.Range("C8").Value = "presencia de horizonte cámbico"

$ws.Range("A9").Value = "cri"
$ws.Range("B9").Value = "frío"
$ws.Range("C9").Value = "frío"

$ws.Range("A10").Value = "dur"
$ws.Range("B10").Value = "duro"
$ws.Range("C10").Value = "con duripán"

$ws.Range("A11").Value = "fibr"
$ws.Range("B11").Value = "fibra"
$ws.Range("C11").Value = "material orgánico poco descompuesto"

$ws.Range("A12").Value = "fluv"
$ws.Range("B12").Value = "rio"
$ws.Range("C12").Value = "llanura de inundación"

$ws.Range("A13").Value = "fol"
$ws.Range("B13").Value = "hojas"
$ws.Range("C13").Value = "masa de hojas"

$ws.Range("A14").Value = "gel"
$ws.Range("B14").Value = "congelado"
$ws.Range("C14").Value = "temperatura media anunal del suelo < 0 ºC"

$ws.Range("A15").Value = "gyps"
$ws.Range("B15").Value = "yeso"
$ws.Range("C15").Value = "presencia de horizonte gypsico"

$ws.Range("A16").Value = "hem"
$ws.Range("B16").Value = "medio"
$ws.Range("C16").Value = "material orgánico semi descompuesto"

$ws.Range("A17").Value = "hist"
$ws.Range("B17").Value = "tejido"
$ws.Range("C17").Value = "presencia de materiales orgánicos"

$ws.Range("A18").Value = "ist"
$ws.Range("B18").Value = "tejido"
$ws.Range("C18").Value = "presencia de materiales orgánicos"

$ws.Range("A19").Value = "hum"
$ws.Range("B19").Value = "humus"
$ws.Range("C19").Value = "presencia de materia orgánica"

$ws.Range("A20").Value = "orth"
$ws.Range("B20").Value = "verdadero"
$ws.Range("C20").Value = "los más comunes"

$ws.Range("A21").Value = "per"
$ws.Range("B21").Value = "todo el año"
$ws.Range("C21").Value = "regimen de humedad perúdico"

$ws.Range("A22").Value = "psamm"
$ws.Range("B22").Value = "arena"
$ws.Range("C22").Value = "texturas arenosas"

$ws.Range("A23").Value = "rend"
$ws.Range("B23").Value = "rendzina, suelo calcáreo"
$ws.Range("C23").Value = "alto contenido de carbonatos"

$ws.Range("A24").Value = "sal"
$ws.Range("B24").Value = "sal"
$ws.Range("C24").Value = "presencia de horizonte sálico"

$ws.Range("A25").Value = "sapr"
$ws.Range("B25").Value = "descompuesto"
$ws.Range("C25").Value = "material orgánico muy descompuesto"

$ws.Range("A26").Value = "torr"
$ws.Range("B26").Value = "cálido, seco"
$ws.Range("C26").Value = "régimen de humedad arídico/tórrico"

$ws.Range("A27").Value = "turb"
$ws.Range("B27").Value = "disturbado"
$ws.Range("C27").Value = "presencia de crioturbación"

$ws.Range("A28").Value = "ud"
$ws.Range("B28").Value = "húmedo"
$ws.Range("C28").Value = "régimen de humedad údico"

$ws.Range("A29").Value = "ust"
$ws.Range("B29").Value = "quemado"
$ws.Range("C29").Value = "régimen de humedad ústico"

$ws.Range("A30").Value = "vitr"
$ws.Range("B30").Value = "vidiro"
$ws.Range("C30").Value = "presencia de vidrio"

$ws.Range("A31").Value = "wass"
$ws.Range("B31").Value = "agua"
$ws.Range("C31").Value = "diariamente bajo aguas poco profundas"

$ws.Range("A32").Value = "xer"
$ws.Range("B32").Value = "seco"
$ws.Range("C32").Value = "régimen de humedad xérico"
